$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row number, then values for columns B,C,D,E ($null means "leave unchanged")
$updates = @(
    @{Row=2;  D="26.918.60"; E="  -0.53%  "},
    @{Row=3;  D="1.668.85";  E="  +0.74%  "},
    @{Row=4;  E="  +0.11%  "},
    @{Row=5;  D="214.64";    E="  -0.32%  "},
    @{Row=6;  D="0.516";     E="  +0.76%  "},
    @{Row=7;  E="  +0.11%  "},
    @{Row=8;  D="0.250";     E="  +0.07%  "},
    @{Row=9;  D="0.0620";    E="  +0.68%  "},
    @{Row=10; D="20.26";     E="  +0.25%  "},
    @{Row=11; E="  +1.51%  "},
    @{Row=12; D="1.905.32";  E="  +0.86%  "},
    @{Row=13; D="1.667.41";  E="  +0.66%  "},
    @{Row=14; E="  -0.04%  "},
    @{Row=15; D="0.528";     E="  +1.46%  "},
    @{Row=16; D="65.52";     E="  +0.34%  "},
    @{Row=17; D="26.939.47"; E="  -0.37%  "},
    @{Row=18; D="233.40";    E="  -1.55%  "},
    @{Row=19; D="7.97";      E="  +2.62%  "},
    @{Row=20; D="0.0₃0732";  E="  +0.20%  "},
    @{Row=21; E="  -0.05%  "},
    @{Row=22; D="4.41";      E="  -0.64%  "},
    @{Row=23; D="9.14";      E="  -2.03%  "},
    @{Row=24; D="2.15";      E="  -2.47%  "},
    @{Row=25; D="146.24";    E="  +0.47%  "},
    @{Row=26; D="7.10";      E="  -0.11%  "},
    @{Row=27; D="15.92";     E="  +0.43%  "},
    @{Row=28; E="  +0.16%  "},
    @{Row=29; D="0.112";     E="  -1.88%  "},
    @{Row=30; D="0.0496";    E="  -0.34%  "},
    @{Row=31; E="  -0.25%  "},
    @{Row=32; E="  -0.02%  "},
    @{Row=33; D="1.450.34";  E="  -7.97%  "},
    @{Row=34; E="  +1.54%  "},
    @{Row=35; E="  +3.12%  "},
    @{Row=36; D="2.41";      E="  +0.03%  "},
    @{Row=37; D="0.585";     E="  +0.92%  "},
    @{Row=38; D="0.897";     E="  -0.51%  "},
    @{Row=39; D="0.0170";    E="  +0.34%  "},
    @{Row=40; D="1.04";      E="  +14.47%  "},
    @{Row=41; D="5.73";      E="  -4.28%  "},
    @{Row=42; E="  +0.07%  "},
    @{Row=43; E="  +1.97%  "},
    @{Row=44; D="66.37";     E="  +1.47%  "},
    @{Row=45; D="1.813.10";  E="  +0.95%  "},
    @{Row=46; D="0.779";     E="  +0.48%  "},
    @{Row=47; D="90.46";     E="  +0.29%  "},
    @{Row=48; E="  +1.26%  "},
    @{Row=49; B="BabyDogeCoin"; C="https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"; D="0.0₆0104"; E="  -0.66%  "},
    @{Row=50; B="Algorand";     C="https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo";     D="0.102";    E="  +3.46%  "},
    @{Row=51; B="Cronos";       C="https://coinranking.com/coin/65PHZTpmE55b+cronos-cro";          D="0.0507";  E="  +0.11%  "}
)

foreach ($u in $updates) {
    $row = $u.Row
    if ($u.ContainsKey("B")) {
        $c = $ws.Cells.Item($row, 2)
        $c.NumberFormat = "@"
        $c.Value = $u.B
    }
    if ($u.ContainsKey("C")) {
        $c = $ws.Cells.Item($row, 3)
        $c.NumberFormat = "@"
        $c.Value = $u.C
    }
    if ($u.ContainsKey("D")) {
        $c = $ws.Cells.Item($row, 4)
        $c.NumberFormat = "@"
        $c.Value = $u.D
    }
    if ($u.ContainsKey("E")) {
        $c = $ws.Cells.Item($row, 5)
        $c.NumberFormat = "@"
        $c.Value = $u.E
    }
}
